$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.102.35"
$ws.Range("E2").Value = "  +0.01%  "

# Row 3
$ws.Range("D3").Value = "1.865.56"
$ws.Range("E3").Value = "  +3.65%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'311.97"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.10%  "

# Row 7
$ws.Range("D7").Value = "'0.4998"
$ws.Range("E7").Value = "  -1.79%  "

# Row 8
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "  +1.26%  "

# Row 9
$ws.Range("D9").Value = "'0.09687"
$ws.Range("E9").Value = "  +25.25%  "

# Row 10
$ws.Range("D10").Value = "'1.128"
$ws.Range("E10").Value = "  +2.89%  "

# Row 11
$ws.Range("D11").Value = "'40.99"
$ws.Range("E11").Value = "  +0.40%  "

# Row 12
$ws.Range("D12").Value = "'6.458"
$ws.Range("E12").Value = "  +2.07%  "

# Row 13
$ws.Range("D13").Value = "'20.90"
$ws.Range("E13").Value = "  +3.31%  "

# Row 14
$ws.Range("D14").Value = "1.865.39"
$ws.Range("E14").Value = "  +3.82%  "

# Row 15
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  -0.05%  "

# Row 16
$ws.Range("D16").Value = "'7.367"
$ws.Range("E16").Value = "  +1.48%  "

# Row 17
$ws.Range("D17").Value = "'0.00001127"
$ws.Range("E17").Value = "  +5.14%  "

# Row 18
$ws.Range("D18").Value = "'92.90"
$ws.Range("E18").Value = "  +1.04%  "

# Row 19
$ws.Range("D19").Value = "'0.06604"
$ws.Range("E19").Value = "  +0.68%  "

# Row 20
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.12%  "

# Row 21
$ws.Range("D21").Value = "'17.39"
$ws.Range("E21").Value = "  +1.19%  "

# Row 22
$ws.Range("D22").Value = "'6.118"
$ws.Range("E22").Value = "  +3.01%  "

# Row 23
$ws.Range("D23").Value = "28.162.64"
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").Value = "'11.26"
$ws.Range("E24").Value = "  +1.97%  "

# Row 25
$ws.Range("D25").Value = "'2.282"
$ws.Range("E25").Value = "  +1.64%  "

# Row 26
$ws.Range("D26").Value = "'2.553"
$ws.Range("E26").Value = "  +5.68%  "

# Row 27
$ws.Range("D27").Value = "2.081.55"
$ws.Range("E27").Value = "  +3.78%  "

# Row 29
$ws.Range("D29").Value = "'157.87"
$ws.Range("E29").Value = "  -1.75%  "

# Row 30
$ws.Range("D30").Value = "'127.13"
$ws.Range("E30").Value = "  +0.07%  "

# Row 31
$ws.Range("D31").Value = "'0.1058"
$ws.Range("E31").Value = "  -2.82%  "

# Row 32
$ws.Range("D32").Value = "'1.059"
$ws.Range("E32").Value = "  +1.53%  "

# Row 33
$ws.Range("D33").Value = "'5.608"
$ws.Range("E33").Value = "  +1.59%  "

# Row 34
$ws.Range("D34").Value = "'3.624"
$ws.Range("E34").Value = "  -0.66%  "

# Row 35
$ws.Range("D35").Value = "'0.06740"
$ws.Range("E35").Value = "  -4.02%  "

# Row 36
$ws.Range("D36").Value = "'9.455"
$ws.Range("E36").Value = "  +4.67%  "

# Row 38
$ws.Range("D38").Value = "'0.2175"
$ws.Range("E38").Value = "  +0.86%  "

# Row 39
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "'11.46"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.996"
$ws.Range("E40").Value = "  -0.37%  "

# Row 41
$ws.Range("D41").Value = "'0.6279"
$ws.Range("E41").Value = "  +3.12%  "

# Row 42
$ws.Range("D42").Value = "'1.174"
$ws.Range("E42").Value = "  +2.03%  "

# Row 43
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44
$ws.Range("D44").Value = "'13.47"
$ws.Range("E44").Value = "  +2.45%  "

# Row 45
$ws.Range("D45").Value = "'0.5985"
$ws.Range("E45").Value = "  +1.60%  "

# Row 46
$ws.Range("D46").Value = "'3.661"
$ws.Range("E46").Value = "  -1.38%  "

# Row 47
$ws.Range("D47").Value = "'1.268"
$ws.Range("E47").Value = "  -1.99%  "

# Row 48
$ws.Range("D48").Value = "'124.24"
$ws.Range("E48").Value = "  -0.60%  "

# Row 49
$ws.Range("D49").Value = "'1.981"
$ws.Range("E49").Value = "  +4.35%  "

# Row 50
$ws.Range("D50").Value = "'1.196"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
$ws.Range("D51").Value = "'0.06828"
$ws.Range("E51").Value = "  +1.32%  "
